# Update the previously-failing test case rows (4, 5, 7, 8) so that their
# Actual Outcome now matches the Expected Outcome ("Same as expected
# outcome.") and the Fail/Pass column reads "Pass" instead of "Fail".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 5, 7, 8)
foreach ($r in $rows) {
    $ws.Range("F$r").Value = "Same as expected outcome."
    $ws.Range("G$r").Value = "Pass"
}

# Move the active selection from F8 to G8.
$ws.Range("G8").Select()
